$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in EARNED value for row 19 (this is reflected in the recalculated
#     BALANCE totals in row 9) ---
$ws.Range("C19").Value = 1.25

# --- Insert a new leave-card entry row above the current row 21, shifting
#     all subsequent PERIOD dates (and the two footer rows at the bottom of
#     the table) down by one row ---
$ws.Rows.Item(21).Insert() | Out-Null

# The freshly inserted row doesn't inherit the surrounding row formatting
# automatically, so copy it over from the row directly below (which now
# holds what used to be row 21's formatting before the insert).
$ws.Range("A22:K22").Copy() | Out-Null
$ws.Range("A21:K21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the calculated-column formula for the new row (table calculated
# columns aren't auto-filled by the insert in this runtime).
$ws.Range("G21").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Grow the table (ListObject) definition to include the newly inserted row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K132")) | Out-Null

# The very last row fell outside the table while it was being resized, so its
# calculated-column formula briefly evaluated outside of table context and
# cached an error; restore it now that the table covers the full range.
$ws.Range("G132").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Fill in the new leave entry's details.
$ws.Range("B21").Value = "FL(2-0-0)"
$ws.Range("D21").Value = 2
$ws.Range("K21").Value = "10/13,16/2023"

# Update the view state to reflect where the user ended up after the edit.
$ws.Range("K21").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
